$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 49.96455633333333
$ws.Range("H2").Value = 149.893669
$ws.Range("I2").Value = 0.1551859508057627
$ws.Range("J2").Value = 0.1551859508057627
$ws.Range("M2").Value = 1.815761
$ws.Range("N2").Value = 5.447283000000001
$ws.Range("O2").Value = 0.07007596730428067
$ws.Range("P2").Value = 0.07007596730428067
$ws.Range("Q2").Value = 90.72369277236966
$ws.Range("R2").Value = 816.513234951327
$ws.Range("S2").Value = 0.01087480561474834
$ws.Range("T2").Value = 0.01087480561474834
$ws.Range("G3").Value = 49.96455633333333
$ws.Range("H3").Value = 149.893669
$ws.Range("I3").Value = 0.1551859508057627
$ws.Range("J3").Value = 0.1551859508057627
$ws.Range("O3").Value = 0.5079540516959071
$ws.Range("P3").Value = 0.5079540516959072
$ws.Range("Q3").Value = 657.6215655852217
$ws.Range("R3").Value = 5918.594090266995
$ws.Range("S3").Value = 0.0788273324780689
$ws.Range("T3").Value = 0.07882733247806892
$ws.Range("G4").Value = 49.96455633333333
$ws.Range("H4").Value = 149.893669
$ws.Range("I4").Value = 0.1551859508057627
$ws.Range("J4").Value = 0.1551859508057627
$ws.Range("M4").Value = 9.711409333333334
$ws.Range("N4").Value = 29.134228
$ws.Range("O4").Value = 0.3747940411327002
$ws.Range("P4").Value = 0.3747940411327002
$ws.Range("Q4").Value = 485.2262587113924
$ws.Range("R4").Value = 4367.036328402532
$ws.Range("S4").Value = 0.05816276962951222
$ws.Range("T4").Value = 0.05816276962951222
$ws.Range("G5").Value = 49.96455633333333
$ws.Range("H5").Value = 149.893669
$ws.Range("I5").Value = 0.1551859508057627
$ws.Range("J5").Value = 0.1551859508057627
$ws.Range("M5").Value = 1.222391
$ws.Range("N5").Value = 3.667173
$ws.Range("O5").Value = 0.04717593986711188
$ws.Range("P5").Value = 0.04717593986711189
$ws.Range("Q5").Value = 61.07622398085966
$ws.Range("R5").Value = 549.6860158277369
$ws.Range("S5").Value = 0.007321043083433245
$ws.Range("T5").Value = 0.007321043083433246
$ws.Range("I6").Value = 0.5303393919600503
$ws.Range("J6").Value = 0.5303393919600503
$ws.Range("M6").Value = 1.815761
$ws.Range("N6").Value = 5.447283000000001
$ws.Range("O6").Value = 0.07007596730428067
$ws.Range("P6").Value = 0.07007596730428067
$ws.Range("Q6").Value = 310.0431953501441
$ws.Range("R6").Value = 2790.388758151297
$ws.Range("S6").Value = 0.03716404589116458
$ws.Range("T6").Value = 0.03716404589116458
$ws.Range("I7").Value = 0.5303393919600503
$ws.Range("J7").Value = 0.5303393919600503
$ws.Range("O7").Value = 0.5079540516959071
$ws.Range("P7").Value = 0.5079540516959072
$ws.Range("S7").Value = 0.2693880429200514
$ws.Range("T7").Value = 0.2693880429200514
$ws.Range("I8").Value = 0.5303393919600503
$ws.Range("J8").Value = 0.5303393919600503
$ws.Range("M8").Value = 9.711409333333334
$ws.Range("N8").Value = 29.134228
$ws.Range("O8").Value = 0.3747940411327002
$ws.Range("P8").Value = 0.3747940411327002
$ws.Range("Q8").Value = 1658.233865062571
$ws.Range("R8").Value = 14924.10478556314
$ws.Range("S8").Value = 0.1987680438845663
$ws.Range("T8").Value = 0.1987680438845663
$ws.Range("I9").Value = 0.5303393919600503
$ws.Range("J9").Value = 0.5303393919600503
$ws.Range("M9").Value = 1.222391
$ws.Range("N9").Value = 3.667173
$ws.Range("O9").Value = 0.04717593986711188
$ws.Range("P9").Value = 0.04717593986711189
$ws.Range("Q9").Value = 208.724612769664
$ws.Range("R9").Value = 1878.521514926976
$ws.Range("S9").Value = 0.02501925926426801
$ws.Range("T9").Value = 0.02501925926426802
$ws.Range("G10").Value = 13.36927633333333
$ws.Range("H10").Value = 40.107829
$ws.Range("I10").Value = 0.04152391238164931
$ws.Range("J10").Value = 0.04152391238164931
$ws.Range("M10").Value = 1.815761
$ws.Range("N10").Value = 5.447283000000001
$ws.Range("O10").Value = 0.07007596730428067
$ws.Range("P10").Value = 0.07007596730428067
$ws.Range("Q10").Value = 24.27541056428966
$ws.Range("R10").Value = 218.478695078607
$ws.Range("S10").Value = 0.002909828326402272
$ws.Range("T10").Value = 0.002909828326402272
$ws.Range("G11").Value = 13.36927633333333
$ws.Range("H11").Value = 40.107829
$ws.Range("I11").Value = 0.04152391238164931
$ws.Range("J11").Value = 0.04152391238164931
$ws.Range("O11").Value = 0.5079540516959071
$ws.Range("P11").Value = 0.5079540516959072
$ws.Range("Q11").Value = 175.9632242987151
$ws.Range("R11").Value = 1583.669018688436
$ws.Range("S11").Value = 0.02109223953652461
$ws.Range("T11").Value = 0.02109223953652462
$ws.Range("G12").Value = 13.36927633333333
$ws.Range("H12").Value = 40.107829
$ws.Range("I12").Value = 0.04152391238164931
$ws.Range("J12").Value = 0.04152391238164931
$ws.Range("M12").Value = 9.711409333333334
$ws.Range("N12").Value = 29.134228
$ws.Range("O12").Value = 0.3747940411327002
$ws.Range("P12").Value = 0.3747940411327002
$ws.Range("Q12").Value = 129.8345149634458
$ws.Range("R12").Value = 1168.510634671012
$ws.Range("S12").Value = 0.01556291492515851
$ws.Range("T12").Value = 0.01556291492515851
$ws.Range("G13").Value = 13.36927633333333
$ws.Range("H13").Value = 40.107829
$ws.Range("I13").Value = 0.04152391238164931
$ws.Range("J13").Value = 0.04152391238164931
$ws.Range("M13").Value = 1.222391
$ws.Range("N13").Value = 3.667173
$ws.Range("O13").Value = 0.04717593986711188
$ws.Range("P13").Value = 0.04717593986711189
$ws.Range("Q13").Value = 16.34248306637966
$ws.Range("R13").Value = 147.082347597417
$ws.Range("S13").Value = 0.00195892959356391
$ws.Range("T13").Value = 0.00195892959356391
$ws.Range("G14").Value = 87.88078300000001
$ws.Range("H14").Value = 263.642349
$ws.Range("I14").Value = 0.2729507448525377
$ws.Range("J14").Value = 0.2729507448525377
$ws.Range("M14").Value = 1.815761
$ws.Range("N14").Value = 5.447283000000001
$ws.Range("O14").Value = 0.07007596730428067
$ws.Range("P14").Value = 0.07007596730428067
$ws.Range("Q14").Value = 159.570498420863
$ws.Range("R14").Value = 1436.134485787767
$ws.Range("S14").Value = 0.01912728747196549
$ws.Range("T14").Value = 0.01912728747196549
$ws.Range("G15").Value = 87.88078300000001
$ws.Range("H15").Value = 263.642349
$ws.Range("I15").Value = 0.2729507448525377
$ws.Range("J15").Value = 0.2729507448525377
$ws.Range("O15").Value = 0.5079540516959071
$ws.Range("P15").Value = 0.5079540516959072
$ws.Range("Q15").Value = 1156.665891632457
$ws.Range("R15").Value = 10409.99302469212
$ws.Range("S15").Value = 0.1386464367612623
$ws.Range("T15").Value = 0.1386464367612623
$ws.Range("G16").Value = 87.88078300000001
$ws.Range("H16").Value = 263.642349
$ws.Range("I16").Value = 0.2729507448525377
$ws.Range("J16").Value = 0.2729507448525377
$ws.Range("M16").Value = 9.711409333333334
$ws.Range("N16").Value = 29.134228
$ws.Range("O16").Value = 0.3747940411327002
$ws.Range("P16").Value = 0.3747940411327002
$ws.Range("Q16").Value = 853.4462562468415
$ws.Range("R16").Value = 7681.016306221572
$ws.Range("S16").Value = 0.1023003126934632
$ws.Range("T16").Value = 0.1023003126934632
$ws.Range("G17").Value = 87.88078300000001
$ws.Range("H17").Value = 263.642349
$ws.Range("I17").Value = 0.2729507448525377
$ws.Range("J17").Value = 0.2729507448525377
$ws.Range("M17").Value = 1.222391
$ws.Range("N17").Value = 3.667173
$ws.Range("O17").Value = 0.04717593986711188
$ws.Range("P17").Value = 0.04717593986711189
$ws.Range("Q17").Value = 107.424678212153
$ws.Range("R17").Value = 966.8221039093771
$ws.Range("S17").Value = 0.01287670792584671
$ws.Range("T17").Value = 0.01287670792584672
